# ElectricityGenerationAndConsumption.xlsx edit
#
# The stray two-column (AW:AX, years 1976/1975) fragment that was sitting
# next to the "Footnotes" block (rows 11-30) is moved back into the main
# data table (rows 34-52), which previously only went back to 1977 (col AV).
# This gives pandas a single contiguous rectangular table to read instead of
# leftover numeric data interleaved with footnote text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Capture the stray AW/AX values from the footnotes area (rows 11-29)
#    before we clear them, so we can re-insert them in the main table.
# ---------------------------------------------------------------------
$strayAW = @{}
$strayAX = @{}
for ($r = 11; $r -le 29; $r++) {
    $strayAW[$r] = $ws.Cells.Item($r, 49).Value2
    $strayAX[$r] = $ws.Cells.Item($r, 50).Value2
}

# ---------------------------------------------------------------------
# 2) Remove the stray AW/AX cells from rows 11-29 entirely (Clear, not
#    ClearContents, so the <c> element itself disappears rather than
#    being left behind as an empty styled cell). Rows 17 and 27 only
#    ever held AW/AX data, so clearing them drops the row completely.
# ---------------------------------------------------------------------
for ($r = 11; $r -le 29; $r++) {
    $ws.Cells.Item($r, 49).Clear()
    $ws.Cells.Item($r, 50).Clear()
}

# Row 11 had an explicit 16pt row height (for the old AW11/AX11 header
# cells); restore it to the sheet's default now that it is a plain
# one-column footnote row again.
$ws.Rows.Item(11).AutoFit()

# ---------------------------------------------------------------------
# 3) Extend the main table (rows 34-52) with two more columns (AW, AX)
#    carrying on from AV (1977) down to 1976 and 1975. Pull the cell
#    formatting from column AV so the new cells match the rest of the
#    table (header style on row 34, data style on rows 35-52).
# ---------------------------------------------------------------------
$ws.Range("AV34:AV52").Copy($ws.Range("AW34:AW52"))
$ws.Range("AV34:AV52").Copy($ws.Range("AX34:AX52"))

# Row 34 is the header row -> "1976 " / "1975 " (note trailing space, to
# match the existing year-header strings). Assign via a literal-text
# formula + paste-values so Excel's numeric autodetection doesn't turn
# "1976 "/"1975 " into the numbers 1976/1975 (which would also mint a
# brand-new, unwanted number-format style instead of reusing the header
# style that was just copied from AV34).
$ws.Range("AW34").Formula = '="1976 "'
$ws.Range("AX34").Formula = '="1975 "'
$ws.Range("AW34:AX34").Copy()
$ws.Range("AW34:AX34").PasteSpecial(-4163)  # xlPasteValues

# Remaining data rows: plain values (numbers / "na" text) straight from
# the captured stray data - none of these collide with the numeric
# auto-detection problem above.
for ($r = 35; $r -le 52; $r++) {
    $srcRow = $r - 23   # 35->12, 36->13, ... 52->29
    $ws.Cells.Item($r, 49).Value = $strayAW[$srcRow]
    $ws.Cells.Item($r, 50).Value = $strayAX[$srcRow]
}

# ---------------------------------------------------------------------
# 4) Update the sheet view to where the edit left the cursor/scroll.
# ---------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 1
$sheetView = $ws.Application.ActiveWindow
$ws.Range("AW55").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AG1").Select()
$ws.Range("AW55").Select()
